$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/61fed3f17293dda63a678ca75783d73e5ed59ff0/e2e/ca88e35a-50f5-4e5b-bf49-6e2dd7945167.md"
$zhcnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc0b96f2e7434503bf88fd73f404366609def715/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ca88e35a-50f5-4e5b-bf49-6e2dd7945167.160a3ba5e00e994264f3ac66c56d0bf80b925c41.zh-cn.xlf"
$dedeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/726097ce4f1001bcfd3d7d687ffa75a7dcdd81f1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ca88e35a-50f5-4e5b-bf49-6e2dd7945167.160a3ba5e00e994264f3ac66c56d0bf80b925c41.de-de.xlf"

$mdDisplay = "ca88e35a-50f5-4e5b-bf49-6e2dd7945167.md"
$zhcnXlfDisplay = "ca88e35a-50f5-4e5b-bf49-6e2dd7945167.160a3ba5e00e994264f3ac66c56d0bf80b925c41.zh-cn.xlf"
$dedeXlfDisplay = "ca88e35a-50f5-4e5b-bf49-6e2dd7945167.160a3ba5e00e994264f3ac66c56d0bf80b925c41.de-de.xlf"

# --- Overview sheet: refresh the "Ready for handoff" status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

# Latest Target File / Latest Handback File columns (E/F) now populated, same
# files as the handoff source (Source File Name / Latest Handoff File)
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl, "", "", $mdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhcnXlfUrl, "", "", $zhcnXlfDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $mdUrl, "", "", $mdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhcnXlfUrl, "", "", $zhcnXlfDisplay)

# Latest Handback DateTime (G) updated to the handback timestamp
$wsZh.Range("G2").Value = "2016-01-26 05:01:26"
$wsZh.Range("G3").Value = "2016-01-26 05:01:26"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl, "", "", $mdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $dedeXlfUrl, "", "", $dedeXlfDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $mdUrl, "", "", $mdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $dedeXlfUrl, "", "", $dedeXlfDisplay)

$wsDe.Range("G2").Value = "2016-01-26 05:01:40"
$wsDe.Range("G3").Value = "2016-01-26 05:01:40"

Write-Output "Report generated for handback"
